# Updates cryptos list values per the Thu Nov  9 21:55:47 UTC 2023 GitHub Actions refresh.
# Rows 27/28 (Cosmos <-> EthereumClassic) also swap position/data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.518.30'
$ws.Range('E2').Value = '  +2.71%  '
$ws.Range('D3').Value = '2.062.94'
$ws.Range('E3').Value = '  +9.29%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''246.56'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').Value = '''0.662'
$ws.Range('E6').Value = '  -4.55%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '''44.90'
$ws.Range('E8').Value = '  +4.31%  '
$ws.Range('D9').Value = '''60.69'
$ws.Range('E9').Value = '  +7.07%  '
$ws.Range('D10').Value = '''0.364'
$ws.Range('E10').Value = '  +2.31%  '
$ws.Range('E11').Value = '  -4.54%  '
$ws.Range('D12').Value = '''0.0984'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '''14.40'
$ws.Range('E13').Value = '  -2.80%  '
$ws.Range('D14').Value = '2.368.05'
$ws.Range('E14').Value = '  +9.56%  '
$ws.Range('D15').Value = '''0.814'
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('D16').Value = '2.056.52'
$ws.Range('E16').Value = '  +9.00%  '
$ws.Range('D17').Value = '''4.89'
$ws.Range('E17').Value = '  -2.84%  '
$ws.Range('D18').Value = '36.544.73'
$ws.Range('E18').Value = '  +2.86%  '
$ws.Range('D19').Value = '''71.20'
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('E20').Value = '  -2.21%  '
$ws.Range('E21').Value = '  -3.52%  '
$ws.Range('E22').Value = '  -3.38%  '
$ws.Range('D23').Value = '''4.88'
$ws.Range('E23').Value = '  -5.73%  '
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('E25').Value = '  -8.10%  '
$ws.Range('D26').Value = '''169.01'
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''20.21'
$ws.Range('E27').Value = '  +10.02%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D28').Value = '''8.81'
$ws.Range('E28').Value = '  +1.84%  '
$ws.Range('D29').Value = '''1.94'
$ws.Range('E29').Value = '  -9.70%  '
$ws.Range('E30').Value = '  -5.22%  '
$ws.Range('D31').Value = '''21.51'
$ws.Range('E31').Value = '  +51.07%  '
$ws.Range('D32').Value = '''4.34'
$ws.Range('E32').Value = '  -1.48%  '
$ws.Range('E33').Value = '  -4.43%  '
$ws.Range('D34').Value = '''0.0906'
$ws.Range('E34').Value = '  +21.49%  '
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').Value = '''2.23'
$ws.Range('E37').Value = '  +14.72%  '
$ws.Range('D38').Value = '''3.98'
$ws.Range('E38').Value = '  -6.78%  '
$ws.Range('D39').Value = '''0.876'
$ws.Range('E39').Value = '  +2.39%  '
$ws.Range('E40').Value = '  -11.19%  '
$ws.Range('E41').Value = '  +2.90%  '
$ws.Range('D42').Value = '''96.58'
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('E43').Value = '  -6.64%  '
$ws.Range('E44').Value = '  +15.84%  '
$ws.Range('D45').Value = '''15.86'
$ws.Range('E45').Value = '  -6.63%  '
$ws.Range('D46').Value = '1.316.03'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('D48').Value = '''2.81'
$ws.Range('E48').Value = '  +2.95%  '
$ws.Range('D49').Value = '2.250.35'
$ws.Range('E49').Value = '  +9.13%  '
$ws.Range('E50').Value = '  -6.47%  '
$ws.Range('E51').Value = '  +14.76%  '
